$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- DOCTOR section: add PHOTO / USER ID columns to the header row ---
$ws.Range("E3").Value = "PHOTO"
$ws.Range("F3").Value = "USER ID"

# --- remove the stray numbered rows (0,1,2) directly under the DOCTOR header ---
$ws.Range("A4").ClearContents()
$ws.Range("A5").ClearContents()
$ws.Range("A6").ClearContents()

# --- TREATMENT section: rework the column headers (row 16) ---
$ws.Range("B16").Value = " TREATMENT TYPE"
$ws.Range("C16").Value = "COST"
$ws.Range("D16").Value = "DOSE"
$ws.Range("E16").Value = "WAY OF ADMINISTRATION"
$ws.Range("F16").Value = "START DATE"
$ws.Range("G16").Value = "END DATE"
$ws.Range("H16").Value = "PATIENT ID"
$ws.Range("I16").Value = "DOCTOR ID"
$ws.Range("J16").Value = "BILL ID"

# --- BILLS section: rework the column headers (row 28) ---
$ws.Range("B28").Value = "TOTAL COST"
$ws.Range("C28").Value = "BANK ID"
$ws.Range("D28").Value = "PAID"
$ws.Range("E28").Value = "PATIENT ID"

# --- NURSE section: add PHOTO / USER ID columns to the header row ---
$ws.Range("E32").Value = "PHOTO"
$ws.Range("F32").Value = "USER ID"

# --- new USER section ---
$ws.Range("A39").Value = "USER"
$ws.Range("A39").Interior.Color = 65535

$ws.Range("A41").Value = "ID"
$ws.Range("B41").Value = "PASSWORD"
$ws.Range("C41").Value = "TYPE"

# --- column widths ---
$ws.Columns.Item(2).ColumnWidth = 13.6
$ws.Columns.Item(4).ColumnWidth = 9.5
$ws.Columns.Item(5).ColumnWidth = 20.1
$ws.Columns.Item(6).ColumnWidth = 20.6

# --- view state: zoom + selection ---
$excel.ActiveWindow.Zoom = 120
$ws.Range("G3").Select() | Out-Null
